$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NetworkLine")

# The R + L + (C // G) branch has been changed to (R + L)//C//G branch.
# Mutual branch rows (10-11): the G column ("inf" meaning no shunt conductance)
# becomes a finite 0 instead of infinite.
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(11, 6).Value = 0

# Self branch rows (12-15): R and wL columns become "inf" (open branch) instead of 0.
$ws.Cells.Item(12, 3).Value = "inf"
$ws.Cells.Item(12, 4).Value = "inf"
$ws.Cells.Item(13, 3).Value = "inf"
$ws.Cells.Item(13, 4).Value = "inf"
$ws.Cells.Item(14, 3).Value = "inf"
$ws.Cells.Item(14, 4).Value = "inf"
$ws.Cells.Item(15, 3).Value = "inf"
$ws.Cells.Item(15, 4).Value = "inf"

# Make NetworkLine the active sheet/tab with F11 selected (matches the saved view state).
$ws.Activate()
$ws.Range("F11").Select()
